$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (D14) and Correspond Handback DateTime (G14)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D14").Value = "2016-03-03 10:52:45"
$wsZhCn.Range("G14").Value = "2016-03-03 10:53:45"

# de-de sheet: update Correspond Handoff Datetime (D14) and Correspond Handback DateTime (G14)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D14").Value = "2016-03-03 10:53:00"
$wsDeDe.Range("G14").Value = "2016-03-03 10:54:09"
